$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 168; this shifts rows 168:235 down to 169:236
$ws.Rows.Item(168).Insert()

# Populate the newly inserted row 168 with the new daily price record
$ws.Range("A168").Value = 5
$ws.Range("B168").Value = "Macroferia Regional de Talca"
$ws.Range("C168").Value = "Maule"
$ws.Range("D168").Value = 45007
$ws.Range("E168").Value = 7
$ws.Range("F168").Value = 100112031
$ws.Range("G168").Value = "Poroto verde"
$ws.Range("H168").Value = "Sin especificar"
$ws.Range("I168").Value = "Primera"
$ws.Range("J168").Value = 100
$ws.Range("K168").Value = 30000
$ws.Range("L168").Value = 30000
$ws.Range("M168").Value = 30000
$ws.Range("N168").Value = "$/saco 25 kilos"
$ws.Range("O168").Value = "Región del Maule"
$ws.Range("P168").Value = 1200
$ws.Range("Q168").Value = 25
$ws.Range("R168").Value = "Hortaliza"
